$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 256, pushing the existing rows 256-267
# down to 258-269 (dimension grows from R267 to R269).
$ws.Range("A256:A257").EntireRow.Insert()

# New row 256: Comercializadora del Agro de Limarí - Ají - Americana (o) - Primera
$ws.Cells.Item(256, 1).Value = 2
$ws.Cells.Item(256, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(256, 3).Value = "Coquimbo"
$ws.Cells.Item(256, 4).Value = 44706
$ws.Cells.Item(256, 5).Value = 4
$ws.Cells.Item(256, 6).Value = 100112021
$ws.Cells.Item(256, 7).Value = "Ají"
$ws.Cells.Item(256, 8).Value = "Americana (o)"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 500
$ws.Cells.Item(256, 11).Value = 24000
$ws.Cells.Item(256, 12).Value = 26000
$ws.Cells.Item(256, 13).Value = 25000
$ws.Cells.Item(256, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(256, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(256, 16).Value = 1000
$ws.Cells.Item(256, 17).Value = 25
$ws.Cells.Item(256, 18).Value = "Hortaliza"

# New row 257: Comercializadora del Agro de Limarí - Ají - Inferno - Primera
$ws.Cells.Item(257, 1).Value = 2
$ws.Cells.Item(257, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(257, 3).Value = "Coquimbo"
$ws.Cells.Item(257, 4).Value = 44706
$ws.Cells.Item(257, 5).Value = 4
$ws.Cells.Item(257, 6).Value = 100112021
$ws.Cells.Item(257, 7).Value = "Ají"
$ws.Cells.Item(257, 8).Value = "Inferno"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 160
$ws.Cells.Item(257, 11).Value = 28000
$ws.Cells.Item(257, 12).Value = 30000
$ws.Cells.Item(257, 13).Value = 29000
$ws.Cells.Item(257, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(257, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(257, 16).Value = 1160
$ws.Cells.Item(257, 17).Value = 25
$ws.Cells.Item(257, 18).Value = "Hortaliza"
